# Generate Report for Archive
#
# 1) The localization status text changed from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:E3 / F2:F3 and the
#    "Status" column (C2:C3) on both the "zh-cn" and "de-de" sheets).
# 2) The two "status" columns got narrower: Overview columns E & F, and
#    column C on the "zh-cn"/"de-de" sheets.

$wb = $excel.ActiveWorkbook

# --- 1) Update the status text on every sheet -----------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2) Narrow the status columns -----------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
